$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the e-mail for Vasily Ivanovich Chapaev (row 5) with a new address,
# which adds a new shared string entry.
$ws.Range("D5").Value = "ivntz.apptest.main@mail.ru"

# Update the active cell selection as recorded in the workbook.
$ws.Range("F13").Select()
